$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes ---

# Row 2 (top-right block)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 8

# Row 6
$ws.Range("I6").Value = 8
$ws.Range("K6").Value = 6

# Row 8
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 6

# Row 12
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 8

# L15, L16, L17 fill color change (red -> orange); new style created first (matches xf index 75)
$ws.Range("L15").Value = 6
$ws.Range("L15").Interior.Color = 49407

$ws.Range("L16").Interior.Color = 49407

$ws.Range("L17").Value = 8
$ws.Range("L17").Interior.Color = 49407

# H15 value + fill color change (orange -> red); new style created next (matches xf index 76)
$ws.Range("H15").Value = 8
$ws.Range("H15").Interior.Color = 255

# H16, H17 value + fill color change (orange -> red); new style created next (matches xf index 77)
$ws.Range("H16").Value = 7
$ws.Range("H16").Interior.Color = 255

$ws.Range("H17").Value = 6
$ws.Range("H17").Interior.Color = 255

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 1
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 6

# --- Selection change ---
$ws.Range("N8").Select()
